$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 2.8
$ws.Range("H3").Value = 3.3
$ws.Range("K3").Value = 2.1
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("AM3").Value = 29

# Row 4 updates
$ws.Range("O4").Value = 1.2
$ws.Range("P4").Value = 4.33
$ws.Range("Q4").Value = 1.7
$ws.Range("R4").Value = 2.1

# Row 8 updates
$ws.Range("G8").Value = 3
$ws.Range("I8").Value = 2.38
$ws.Range("J8").Value = 3.75
$ws.Range("L8").Value = 3.1
$ws.Range("W8").Value = 8.5
$ws.Range("X8").Value = 15
$ws.Range("AA8").Value = 26
$ws.Range("AI8").Value = 11
$ws.Range("AJ8").Value = 9.5
$ws.Range("AN8").Value = 5
$ws.Range("AW8").Value = 4.33
$ws.Range("AX8").Value = 13
$ws.Range("AZ8").Value = 41
